$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns B..K
$ws.Range("B1").Value = "id"
$ws.Range("C1").Value = "genbank"
$ws.Range("D1").Value = "length"
$ws.Range("E1").Value = "annotated"
$ws.Range("F1").Value = "seedOrtholog"
$ws.Range("G1").Value = "evalueEggnog"
$ws.Range("H1").Value = "scoreEggnog"
$ws.Range("I1").Value = "description"
$ws.Range("J1").Value = "preferredName"
$ws.Range("K1").Value = "clusterId_id"

# Copy B1's formatting (bold, border, alignment) onto the rest of the header row
$ws.Range("B1").Copy()
$ws.Range("C1:K1").PasteSpecial(-4122)

# Data row (row 2)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "QWERTY"
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "ACGATGCTAGTATCG"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = "Sed ut perspiciatis unde omnis iste natus error si"
$ws.Range("J2").Value = "PAR3"
$ws.Range("K2").Value = 1
